$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param(
        [string]$CellRef,
        [string]$Text
    )
    $cell = $ws.Range($CellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $Text
}

# Row 2 - Bitcoin
Set-CellText "D2" "69.323.22"
Set-CellText "E2" "  -2.60%  "

# Row 3 - Ethereum
Set-CellText "D3" "3.526.05"
Set-CellText "E3" "  -4.64%  "

# Row 4 - TetherUSD
Set-CellText "E4" "  +0.04%  "

# Row 5 - BNB
Set-CellText "D5" "579.93"
Set-CellText "E5" "  -0.42%  "

# Row 6 - Solana
Set-CellText "D6" "171.62"
Set-CellText "E6" "  -3.74%  "

# Row 7 - now LidoStakedEther (was XRP)
Set-CellText "B7" "LidoStakedEther"
Set-CellText "C7" "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
Set-CellText "D7" "3.517.65"
Set-CellText "E7" "  -4.60%  "

# Row 8 - now XRP (was LidoStakedEther)
Set-CellText "B8" "XRP"
Set-CellText "C8" "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
Set-CellText "D8" "0.609"
Set-CellText "E8" "  -1.28%  "

# Row 9 - USDC
Set-CellText "E9" "  +0.06%  "

# Row 10 - Dogecoin
Set-CellText "D10" "0.190"
Set-CellText "E10" "  -4.88%  "

# Row 11 - Toncoin
Set-CellText "D11" "6.74"
Set-CellText "E11" "  -3.19%  "

# Row 12 - Cardano
Set-CellText "D12" "0.586"
Set-CellText "E12" "  -4.24%  "

# Row 13 - Avalanche
Set-CellText "D13" "47.37"
Set-CellText "E13" "  -3.71%  "

# Row 14 - ShibaInu
Set-CellText "E14" "  -4.70%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-CellText "D15" "4.091.82"
Set-CellText "E15" "  -4.61%  "

# Row 16 - Polkadot
Set-CellText "D16" "8.55"
Set-CellText "E16" "  -5.24%  "

# Row 17 - BitcoinCash
Set-CellText "D17" "631.26"
Set-CellText "E17" "  -7.10%  "

# Row 18 - WrappedEther
Set-CellText "D18" "3.523.89"
Set-CellText "E18" "  -4.87%  "

# Row 19 - WrappedBTC
Set-CellText "D19" "69.308.85"
Set-CellText "E19" "  -2.77%  "

# Row 20 - TRON
Set-CellText "E20" "  +0.09%  "

# Row 21 - Chainlink
Set-CellText "D21" "17.51"
Set-CellText "E21" "  -2.69%  "

# Row 22 - Uniswap
Set-CellText "D22" "11.22"
Set-CellText "E22" "  -3.37%  "

# Row 23 - Polygon
Set-CellText "D23" "0.888"
Set-CellText "E23" "  -6.01%  "

# Row 24 - InternetComputer(DFINITY)
Set-CellText "D24" "15.99"
Set-CellText "E24" "  -8.28%  "

# Row 25 - Litecoin
Set-CellText "D25" "97.88"
Set-CellText "E25" "  -4.31%  "

# Row 26 - PancakeSwap
Set-CellText "D26" "3.82"
Set-CellText "E26" "  -4.26%  "

# Row 27 - Dai
Set-CellText "E27" "  +0.13%  "

# Row 28 - ImmutableX
Set-CellText "E28" "  -7.13%  "

# Row 29 - RenderToken
Set-CellText "D29" "9.34"
Set-CellText "E29" "  -9.29%  "

# Row 30 - EthereumClassic
Set-CellText "D30" "32.85"
Set-CellText "E30" "  -6.49%  "

# Row 31 - Stacks
Set-CellText "D31" "3.17"
Set-CellText "E31" "  -7.55%  "

# Row 32 - Filecoin
Set-CellText "D32" "8.58"
Set-CellText "E32" "  -6.51%  "

# Row 33 - Mantle
Set-CellText "E33" "  -7.21%  "

# Row 34 - NEARProtocol
Set-CellText "E34" "  -6.81%  "

# Row 35 - Bittensor
Set-CellText "D35" "633.03"
Set-CellText "E35" "  +8.94%  "

# Row 36 - Cosmos
Set-CellText "D36" "10.78"
Set-CellText "E36" "  -3.84%  "

# Row 37 - now dogwifhat (was Hedera)
Set-CellText "B37" "dogwifhat"
Set-CellText "C37" "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-CellText "D37" "3.50"
Set-CellText "E37" "  -14.42%  "

# Row 38 - now Hedera (was dogwifhat)
Set-CellText "B38" "Hedera"
Set-CellText "C38" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-CellText "D38" "0.103"
Set-CellText "E38" "  -4.92%  "

# Row 39 - OKB
Set-CellText "D39" "57.24"
Set-CellText "E39" "  -2.48%  "

# Row 40 - FirstDigitalUSD
Set-CellText "E40" "  +0.11%  "

# Row 41 - VeChain
Set-CellText "D41" "0.0457"
Set-CellText "E41" "  -0.21%  "

# Row 42 - Kaspa
Set-CellText "E42" "  -5.72%  "

# Row 43 - Maker
Set-CellText "D43" "3.385.56"
Set-CellText "E43" "  -8.10%  "

# Row 44 - TheGraph
Set-CellText "D44" "0.330"
Set-CellText "E44" "  -6.76%  "

# Row 45 - InjectiveProtocol
Set-CellText "D45" "33.03"
Set-CellText "E45" "  -7.65%  "

# Row 46 - PEPE
Set-CellText "D46" "0.0$([char]0x2083)0699"
Set-CellText "E46" "  -9.10%  "

# Row 47 - Fetch.AI
Set-CellText "E47" "  -7.59%  "

# Row 48 - ThetaToken
Set-CellText "E48" "  -5.80%  "

# Row 50 - MXToken
Set-CellText "E50" "  +14.83%  "

# Row 51 - Monero
Set-CellText "D51" "131.79"
Set-CellText "E51" "  -2.48%  "
